$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("J5").Value = -0.4662
$ws.Range("K5").Value = -0.0489
$ws.Range("L5").Value = -0.0509
$ws.Range("M5").Value = -0.0429
$ws.Range("N5").Value = -0.295
$ws.Range("O5").Value = -0.4274
$ws.Range("P5").Value = -0.0613
$ws.Range("Q5").Value = -0.0499
$ws.Range("R5").Value = -0.3899

# Row 6
$ws.Range("J6").Value = -0.0641
$ws.Range("K6").Value = -0.0694
$ws.Range("L6").Value = -0.0526
$ws.Range("M6").Value = -0.0471
$ws.Range("N6").Value = -0.0352
$ws.Range("O6").Value = -0.0202
$ws.Range("P6").Value = -0.0223
$ws.Range("Q6").Value = -0.0112
$ws.Range("R6").Value = -0.0101

# Row 7
$ws.Range("J7").Value = -0.1115
$ws.Range("K7").Value = -0.0957
$ws.Range("L7").Value = -0.1426
$ws.Range("M7").Value = -0.1501
$ws.Range("N7").Value = -0.1383
$ws.Range("O7").Value = -0.038
$ws.Range("P7").Value = 0.0151
$ws.Range("Q7").Value = 0.0197
$ws.Range("R7").Value = 0.03

# Row 8
$ws.Range("J8").Value = -0.117
$ws.Range("K8").Value = 0.0557
$ws.Range("L8").Value = 0.1076
$ws.Range("M8").Value = 0.1645
$ws.Range("N8").Value = 0.4363
$ws.Range("O8").Value = 0.4298
$ws.Range("P8").Value = 0.4413
$ws.Range("Q8").Value = 0.4039
$ws.Range("R8").Value = 0.2893

# Row 16
$ws.Range("J16").Value = -1.5464
$ws.Range("K16").Value = -1.6509
$ws.Range("L16").Value = -0.585
$ws.Range("M16").Value = -0.5631
$ws.Range("N16").Value = -0.6118
$ws.Range("O16").Value = -0.1892
$ws.Range("P16").Value = -0.1336
$ws.Range("Q16").Value = 0.0243
$ws.Range("R16").Value = -0.0478

# Row 33
$ws.Range("J33").Value = -0.0693
$ws.Range("K33").Value = 0.0125
$ws.Range("L33").Value = 0.114
$ws.Range("M33").Value = 0.081
$ws.Range("N33").Value = 0.0196
$ws.Range("O33").Value = 0.015
$ws.Range("P33").Value = 0.0236
$ws.Range("Q33").Value = -0.0244
$ws.Range("R33").Value = 0.0356

# Row 34
$ws.Range("J34").Value = 0.04
$ws.Range("K34").Value = 0.0413
$ws.Range("L34").Value = 0.0429
$ws.Range("M34").Value = 0.0448
$ws.Range("N34").Value = 0.0464
$ws.Range("O34").Value = 0.0481
$ws.Range("P34").Value = 0.0495
$ws.Range("Q34").Value = 0.0412
$ws.Range("R34").Value = 0.0365

# Row 35
$ws.Range("J35").Value = -0.0452
$ws.Range("K35").Value = -0.0879
$ws.Range("L35").Value = -0.1301
$ws.Range("M35").Value = -0.1319
$ws.Range("N35").Value = -0.0867
$ws.Range("O35").Value = -0.0414
$ws.Range("P35").Value = 0.0026
$ws.Range("Q35").Value = 0.0076
$ws.Range("R35").Value = 0.0087

# Row 36
$ws.Range("J36").Value = 0.0679
$ws.Range("K36").Value = 0.1303
$ws.Range("L36").Value = 0.1576
$ws.Range("M36").Value = 0.1184
$ws.Range("N36").Value = 0.0922
$ws.Range("O36").Value = 0.0924
$ws.Range("P36").Value = 0.0869
$ws.Range("Q36").Value = 0.061
$ws.Range("R36").Value = -0.0027

# Row 44
$ws.Range("J44").Value = -0.0746
$ws.Range("K44").Value = 0.1391
$ws.Range("L44").Value = 0.2986
$ws.Range("M44").Value = 0.1398
$ws.Range("N44").Value = 0.1308
$ws.Range("O44").Value = 0.1549
$ws.Range("P44").Value = 0.2022
$ws.Range("Q44").Value = 0.1181
$ws.Range("R44").Value = -0.082
